# Apply the cryptos-list price/volume refresh described by the commit diff.
# Only the D (Price) and E (Volume 1h) columns change, rows 2-51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.894.94"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").Value = "3.418.12"
$ws.Range("E3").Value = "  +0.86%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.36"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.71"
$ws.Range("E6").Value = "  +1.90%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.474"
$ws.Range("E8").Value = "  -0.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.61"
$ws.Range("E9").Value = "  -0.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.123"
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.386"
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("D12").Value = "4.006.59"
$ws.Range("E12").Value = "  +1.02%  "
$ws.Range("E13").Value = "  -0.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.21"
$ws.Range("E14").Value = "  +1.19%  "
$ws.Range("D15").Value = "3.419.98"
$ws.Range("E15").Value = "  +0.50%  "
$ws.Range("E16").Value = "  -0.64%  "
$ws.Range("D17").Value = "61.952.41"
$ws.Range("E17").Value = "  +1.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.19"
$ws.Range("E18").Value = "  +1.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.95"
$ws.Range("E19").Value = "  +1.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.19"
$ws.Range("E20").Value = "  +3.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "389.56"
$ws.Range("E21").Value = "  +1.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.30"
$ws.Range("E22").Value = "  -1.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.552"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000115"
$ws.Range("E25").Value = "  -0.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.191"
$ws.Range("E26").Value = "  +3.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.45"
$ws.Range("E27").Value = "  +2.78%  "
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.04"
$ws.Range("E29").Value = "  +0.62%  "
$ws.Range("E30").Value = "  +0.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.42"
$ws.Range("E31").Value = "  +2.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.52"
$ws.Range("E33").Value = "  +1.07%  "
$ws.Range("E34").Value = "  +5.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.97"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "167.92"
$ws.Range("E36").Value = "  +0.91%  "
$ws.Range("D37").Value = "3.454.00"
$ws.Range("E37").Value = "  +1.03%  "
$ws.Range("E38").Value = "  +0.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "28.50"
$ws.Range("E39").Value = "  +7.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0752"
$ws.Range("E40").Value = "  -2.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.786"
$ws.Range("E41").Value = "  +0.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.44"
$ws.Range("E42").Value = "  +1.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.67"
$ws.Range("E43").Value = "  +0.80%  "
$ws.Range("E44").Value = "  +4.23%  "
$ws.Range("D45").Value = "2.532.35"
$ws.Range("E45").Value = "  +3.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.82"
$ws.Range("E46").Value = "  -0.85%  "
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.61"
$ws.Range("E48").Value = "  -1.05%  "
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("E50").Value = "  -3.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.205"
$ws.Range("E51").Value = "  -0.58%  "
